# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F2").Value = 5336
$wsExh.Range("F4").Value = 11106
$wsExh.Range("F7").Value = 158
$wsExh.Range("F8").Value = 228
$wsExh.Range("F9").Value = 952

# Sheet "演出" (performances)
$wsPerf = $wb.Worksheets.Item("演出")
$wsPerf.Range("F2").Value = 16

# Sheet "全部类型" (all types, combined view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 16
$wsAll.Range("F4").Value = 5336
$wsAll.Range("F7").Value = 11106
$wsAll.Range("F10").Value = 158
$wsAll.Range("F13").Value = 228
$wsAll.Range("F14").Value = 952
